# Updated cryptos list with GitHub Actions.
# Price (column D) values are digit-grouped strings (e.g. "27.912.43") that must stay
# as text rather than being auto-parsed as numbers by Excel, so a leading apostrophe
# (PowerShell-escaped as '' inside a single-quoted string) is used to force text entry
# for those cells, matching the original inline-string cell content.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.912.43'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '''1.811.23'
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("D4").Value = '''1.001'
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").Value = '''309.62'
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").Value = '''0.9997'
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").Value = '''0.4939'
$ws.Range("E7").Value = '  -3.49%  '
$ws.Range("D8").Value = '''0.3869'
$ws.Range("E8").Value = '  +2.75%  '
$ws.Range("D9").Value = '''0.09805'
$ws.Range("E9").Value = '  +26.08%  '
$ws.Range("E10").Value = '  +1.49%  '
$ws.Range("D11").Value = '''40.86'
$ws.Range("D12").Value = '''6.434'
$ws.Range("E12").Value = '  +4.01%  '
$ws.Range("D13").Value = '''20.50'
$ws.Range("E13").Value = '  +1.65%  '
$ws.Range("D14").Value = '''0.9993'
$ws.Range("E14").Value = '  -0.40%  '
$ws.Range("D15").Value = '''1.810.94'
$ws.Range("E15").Value = '  +2.11%  '
$ws.Range("D16").Value = '''7.288'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("D17").Value = '''0.00001135'
$ws.Range("E17").Value = '  +6.21%  '
$ws.Range("D18").Value = '''92.59'
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").Value = '''0.06602'
$ws.Range("E19").Value = '  +1.08%  '
$ws.Range("D20").Value = '''0.9996'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").Value = '''5.934'
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("D23").Value = '''27.973.74'
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").Value = '''11.13'
$ws.Range("E24").Value = '  +1.77%  '
$ws.Range("D25").Value = '''2.241'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("D26").Value = '''158.78'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").Value = '''2.020.16'
$ws.Range("E27").Value = '  +2.02%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").Value = '''20.56'
$ws.Range("E28").Value = '  +1.61%  '
$ws.Range("D29").Value = '''2.390'
$ws.Range("E29").Value = '  +1.74%  '
$ws.Range("D30").Value = '''127.27'
$ws.Range("E30").Value = '  +4.12%  '
$ws.Range("D31").Value = '''0.1059'
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").Value = '''5.577'
$ws.Range("E33").Value = '  +1.77%  '
$ws.Range("D34").Value = '''3.630'
$ws.Range("E34").Value = '  -0.03%  '
$ws.Range("D35").Value = '''0.06756'
$ws.Range("E35").Value = '  -4.66%  '
$ws.Range("D36").Value = '''9.021'
$ws.Range("E36").Value = '  +5.20%  '
$ws.Range("D37").Value = '''0.02325'
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("D38").Value = '''0.2132'
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").Value = '''4.938'
$ws.Range("E39").Value = '  -1.36%  '
$ws.Range("D40").Value = '''11.26'
$ws.Range("E40").Value = '  -2.05%  '
$ws.Range("D41").Value = '''0.6198'
$ws.Range("E41").Value = '  +1.65%  '
$ws.Range("D42").Value = '''0.9994'
$ws.Range("E42").Value = '  -0.39%  '
$ws.Range("D43").Value = '''1.142'
$ws.Range("E43").Value = '  -0.92%  '
$ws.Range("D44").Value = '''13.04'
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '''0.5865'
$ws.Range("E45").Value = '  -1.65%  '
$ws.Range("D46").Value = '''3.688'
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("D47").Value = '''1.278'
$ws.Range("E47").Value = '  -4.00%  '
$ws.Range("D48").Value = '''122.36'
$ws.Range("E48").Value = '  -3.15%  '
$ws.Range("D49").Value = '''1.931'
$ws.Range("E49").Value = '  +1.97%  '
$ws.Range("D50").Value = '''1.174'
$ws.Range("E50").Value = '  -3.24%  '
$ws.Range("E51").Value = '  +1.17%  '
